# Updates market-price-derived columns (H-N) across several sheets of the
# "Lich_Profits" workbook, as refreshed by the scheduled market-data runner.
# Columns: H=currentAveragePrice, I=currentAveragePriceNQ, J=currentAveragePriceHQ,
#          K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ, N=LeveProfitHQ

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 439.5
$ws.Range("I2").Value = 399.5
$ws.Range("J2").Value = 479.5
$ws.Range("K2").Value = 399.5
$ws.Range("L2").Value = 479.5
$ws.Range("M2").Value = -286.5
$ws.Range("N2").Value = -705.5

$ws.Range("H70").Value = 8996.25
$ws.Range("J70").Value = 9981.429
$ws.Range("L70").Value = 29944.287
$ws.Range("N70").Value = -30484.287

$ws.Range("H73").Value = 8996.25
$ws.Range("J73").Value = 9981.429
$ws.Range("L73").Value = 29944.287
$ws.Range("N73").Value = -31816.287

$ws.Range("H112").Value = 4102247.8
$ws.Range("J112").Value = 4102247.8
$ws.Range("L112").Value = 12306743.4
$ws.Range("N112").Value = -12308959.4

$ws.Range("H113").Value = 7629.087
$ws.Range("I113").Value = 12303.583
$ws.Range("J113").Value = 2529.6365
$ws.Range("K113").Value = 12303.583
$ws.Range("L113").Value = 2529.6365
$ws.Range("M113").Value = -9049.583000000001
$ws.Range("N113").Value = -9037.636500000001

$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").Value = $null

$ws.Range("H137").Value = 50217.39
$ws.Range("I137").Value = 93376.5
$ws.Range("J137").Value = 3134.7273
$ws.Range("K137").Value = 280129.5
$ws.Range("L137").Value = 9404.1819
$ws.Range("M137").Value = -277579.5
$ws.Range("N137").Value = -14504.1819

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 8627.666999999999
$ws.Range("I45").Value = 15048.625
$ws.Range("K45").Value = 15048.625
$ws.Range("M45").Value = -14671.625

$ws.Range("H61").Value = 3062.9565
$ws.Range("I61").Value = 2092
$ws.Range("J61").Value = 7675
$ws.Range("K61").Value = 2092
$ws.Range("L61").Value = 7675
$ws.Range("M61").Value = -1880
$ws.Range("N61").Value = -8099

$ws.Range("H110").Value = 3458.5144
$ws.Range("I110").Value = 3470.1924
$ws.Range("J110").Value = 3424.7778
$ws.Range("K110").Value = 3470.1924
$ws.Range("L110").Value = 3424.7778
$ws.Range("M110").Value = -1425.1924
$ws.Range("N110").Value = -7514.7778

$ws.Range("H122").Value = 7038.174
$ws.Range("I122").Value = 6914.5
$ws.Range("K122").Value = 20743.5
$ws.Range("M122").Value = -18293.5

$ws.Range("H132").Value = 4255
$ws.Range("I132").Value = 4107.6
$ws.Range("K132").Value = 12322.8
$ws.Range("M132").Value = -9792.800000000001

$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").Value = $null

$ws.Range("H136").Value = 3062.9565
$ws.Range("I136").Value = 2092
$ws.Range("J136").Value = 7675
$ws.Range("K136").Value = 6276
$ws.Range("L136").Value = 23025
$ws.Range("M136").Value = -3726
$ws.Range("N136").Value = -28125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2797.5527
$ws.Range("I86").Value = 2428.2173
$ws.Range("J86").Value = 3363.8667
$ws.Range("K86").Value = 2428.2173
$ws.Range("L86").Value = 3363.8667
$ws.Range("M86").Value = -1305.2173
$ws.Range("N86").Value = -5609.8667

$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").Value = $null

$ws.Range("H89").Value = 2797.5527
$ws.Range("I89").Value = 2428.2173
$ws.Range("J89").Value = 3363.8667
$ws.Range("K89").Value = 12141.0865
$ws.Range("L89").Value = 16819.3335
$ws.Range("M89").Value = -6525.086499999999
$ws.Range("N89").Value = -28051.3335

$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").Value = $null

$ws.Range("H135").Value = 66666.336
$ws.Range("J135").Value = 66666.336
$ws.Range("L135").Value = 66666.336
$ws.Range("N135").Value = -76806.336

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 146965.88
$ws.Range("I31").Value = 335097.44
$ws.Range("J31").Value = 2249.282
$ws.Range("K31").Value = 335097.44
$ws.Range("L31").Value = 2249.282
$ws.Range("M31").Value = -334802.44
$ws.Range("N31").Value = -2839.282

$ws.Range("H34").Value = 146965.88
$ws.Range("I34").Value = 335097.44
$ws.Range("J34").Value = 2249.282
$ws.Range("K34").Value = 335097.44
$ws.Range("L34").Value = 2249.282
$ws.Range("M34").Value = -334895.44
$ws.Range("N34").Value = -2653.282

$ws.Range("H58").Value = 2688.8572
$ws.Range("I58").Value = 2364.6
$ws.Range("K58").Value = 2364.6
$ws.Range("M58").Value = -2161.6

$ws.Range("H105").Value = 3835.9048
$ws.Range("I105").Value = 1800.8572
$ws.Range("K105").Value = 1800.8572
$ws.Range("M105").Value = -53.85719999999992

$ws.Range("H134").Value = 3394.4583
$ws.Range("I134").Value = 2940.6316
$ws.Range("J134").Value = 5119
$ws.Range("K134").Value = 8821.8948
$ws.Range("L134").Value = 15357
$ws.Range("M134").Value = -6286.8948
$ws.Range("N134").Value = -20427

$ws.Range("H136").Value = 2688.8572
$ws.Range("I136").Value = 2364.6
$ws.Range("K136").Value = 7093.799999999999
$ws.Range("M136").Value = -4543.799999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 113.333336
$ws.Range("I40").Value = 111
$ws.Range("K40").Value = 444
$ws.Range("M40").Value = -375

$ws.Range("H58").Value = 1333.3334
$ws.Range("I58").Value = 500
$ws.Range("J58").Value = 1500
$ws.Range("K58").Value = 1500
$ws.Range("L58").Value = 4500
$ws.Range("M58").Value = -1372
$ws.Range("N58").Value = -4756

$ws.Range("H131").Value = 11906352
$ws.Range("J131").Value = 1771.7931
$ws.Range("L131").Value = 5315.379300000001
$ws.Range("N131").Value = -15395.3793

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 29999.75
$ws.Range("I57").Value = 20000
$ws.Range("K57").Value = 20000
$ws.Range("M57").Value = -19180

$ws.Range("H134").Value = 40121
$ws.Range("J134").Value = 40121
$ws.Range("L134").Value = 120363
$ws.Range("N134").Value = -125433

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5397.8
$ws.Range("I122").Value = 4248.5
$ws.Range("K122").Value = 12745.5
$ws.Range("M122").Value = -10295.5

$ws.Range("H132").Value = 3030.9583
$ws.Range("I132").Value = 2692.325
$ws.Range("K132").Value = 8076.974999999999
$ws.Range("M132").Value = -5546.974999999999

